# Auto-generated Excel COM-interop script to refresh market/profit data cells
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets (WVR unchanged).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1626.75
$ws.Range("I6").Value = 1502.3334
$ws.Range("K6").Value = 4507.0002
$ws.Range("M6").Value = -4395.0002
$ws.Range("H8").Value = 3241.7144
$ws.Range("I8").Value = 1930.1666
$ws.Range("K8").Value = 5790.4998
$ws.Range("M8").Value = -5651.4998
$ws.Range("H43").Value = 1040
$ws.Range("I43").Value = 1050
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 1050
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -981
$ws.Range("N43").Value = -1138
$ws.Range("H62").Value = 2863956
$ws.Range("I62").Value = 4127473.5
$ws.Range("J62").Value = 21041.584
$ws.Range("K62").Value = 4127473.5
$ws.Range("L62").Value = 21041.584
$ws.Range("M62").Value = -4126849.5
$ws.Range("N62").Value = -22289.584
$ws.Range("H65").Value = 2863956
$ws.Range("I65").Value = 4127473.5
$ws.Range("J65").Value = 21041.584
$ws.Range("K65").Value = 20637367.5
$ws.Range("L65").Value = 105207.92
$ws.Range("M65").Value = -20634247.5
$ws.Range("N65").Value = -111447.92
$ws.Range("H113").Value = 5798.273
$ws.Range("I113").Value = 4518.3335
$ws.Range("J113").Value = 6278.25
$ws.Range("K113").Value = 4518.3335
$ws.Range("L113").Value = 6278.25
$ws.Range("M113").Value = -1264.3335
$ws.Range("N113").Value = -12786.25
$ws.Range("H129").Value = 1206.6957
$ws.Range("I129").Value = 290.66666
$ws.Range("J129").Value = 1344.1
$ws.Range("K129").Value = 871.9999799999999
$ws.Range("L129").Value = 4032.3
$ws.Range("M129").Value = 4128.00002
$ws.Range("N129").Value = -14032.3
$ws.Range("H132").Value = 24789.232
$ws.Range("I132").Value = 27011.205
$ws.Range("J132").Value = 3125
$ws.Range("K132").Value = 81033.61500000001
$ws.Range("L132").Value = 9375
$ws.Range("M132").Value = -78503.61500000001
$ws.Range("N132").Value = -14435

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33463.13
$ws.Range("I2").Value = 41294.72
$ws.Range("J2").Value = 831.5
$ws.Range("K2").Value = 41294.72
$ws.Range("L2").Value = 831.5
$ws.Range("M2").Value = -41181.72
$ws.Range("N2").Value = -1057.5
$ws.Range("H61").Value = 3261
$ws.Range("I61").Value = 1881.2
$ws.Range("K61").Value = 1881.2
$ws.Range("M61").Value = -1669.2
$ws.Range("H74").Value = 4138.7383
$ws.Range("I74").Value = 1084.1111
$ws.Range("K74").Value = 1084.1111
$ws.Range("M74").Value = -210.1111000000001
$ws.Range("H77").Value = 4138.7383
$ws.Range("I77").Value = 1084.1111
$ws.Range("K77").Value = 5420.5555
$ws.Range("M77").Value = -1052.5555
$ws.Range("H116").Value = 33463.13
$ws.Range("I116").Value = 41294.72
$ws.Range("J116").Value = 831.5
$ws.Range("K116").Value = 41294.72
$ws.Range("L116").Value = 831.5
$ws.Range("M116").Value = -39000.72
$ws.Range("N116").Value = -5419.5
$ws.Range("H122").Value = 2135.7646
$ws.Range("I122").Value = 2214.0715
$ws.Range("J122").Value = 1770.3334
$ws.Range("K122").Value = 6642.2145
$ws.Range("L122").Value = 5311.0002
$ws.Range("M122").Value = -4192.2145
$ws.Range("N122").Value = -10211.0002
$ws.Range("H136").Value = 3261
$ws.Range("I136").Value = 1881.2
$ws.Range("K136").Value = 5643.6
$ws.Range("M136").Value = -3093.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33463.13
$ws.Range("I3").Value = 41294.72
$ws.Range("J3").Value = 831.5
$ws.Range("K3").Value = 41294.72
$ws.Range("L3").Value = 831.5
$ws.Range("M3").Value = -41180.72
$ws.Range("N3").Value = -1059.5
$ws.Range("H20").Value = 1479.7333
$ws.Range("I20").Value = 1319.6
$ws.Range("K20").Value = 1319.6
$ws.Range("M20").Value = -1072.6
$ws.Range("H86").Value = 1567.1578
$ws.Range("I86").Value = 1586.5555
$ws.Range("J86").Value = 1549.7
$ws.Range("K86").Value = 1586.5555
$ws.Range("L86").Value = 1549.7
$ws.Range("M86").Value = -463.5554999999999
$ws.Range("N86").Value = -3795.7
$ws.Range("H89").Value = 1567.1578
$ws.Range("I89").Value = 1586.5555
$ws.Range("J89").Value = 1549.7
$ws.Range("K89").Value = 7932.7775
$ws.Range("L89").Value = 7748.5
$ws.Range("M89").Value = -2316.7775
$ws.Range("N89").Value = -18980.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4769.1787
$ws.Range("I31").Value = 1511.5454
$ws.Range("J31").Value = 9443.174000000001
$ws.Range("K31").Value = 1511.5454
$ws.Range("L31").Value = 9443.174000000001
$ws.Range("M31").Value = -1216.5454
$ws.Range("N31").Value = -10033.174
$ws.Range("H34").Value = 4769.1787
$ws.Range("I34").Value = 1511.5454
$ws.Range("J34").Value = 9443.174000000001
$ws.Range("K34").Value = 1511.5454
$ws.Range("L34").Value = 9443.174000000001
$ws.Range("M34").Value = -1309.5454
$ws.Range("N34").Value = -9847.174000000001
$ws.Range("H99").Value = 2878.5
$ws.Range("I99").Value = 1928.4546
$ws.Range("J99").Value = 4371.4287
$ws.Range("K99").Value = 1928.4546
$ws.Range("L99").Value = 4371.4287
$ws.Range("M99").Value = -430.4546
$ws.Range("N99").Value = -7367.4287
$ws.Range("H122").Value = 1804
$ws.Range("I122").Value = 1886.4
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 5659.200000000001
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -3209.200000000001
$ws.Range("N122").Value = -9899.9998
$ws.Range("H126").Value = 2878.5
$ws.Range("I126").Value = 1928.4546
$ws.Range("J126").Value = 4371.4287
$ws.Range("K126").Value = 5785.3638
$ws.Range("L126").Value = 13114.2861
$ws.Range("M126").Value = -3315.3638
$ws.Range("N126").Value = -18054.2861
$ws.Range("H134").Value = 1995.3462
$ws.Range("I134").Value = 1269.5227
$ws.Range("K134").Value = 3808.5681
$ws.Range("M134").Value = -1273.5681

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1978.75
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1978.75
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 5936.25
$ws.Range("N75").Value = -7932.25
$ws.Range("H78").Value = 1978.75
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1978.75
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 17808.75
$ws.Range("N78").Value = -27792.75
$ws.Range("H103").Value = 232257.92
$ws.Range("I103").Value = 1048.7
$ws.Range("J103").Value = 1002955.3
$ws.Range("K103").Value = 3146.1
$ws.Range("L103").Value = 3008865.9
$ws.Range("M103").Value = -2267.1
$ws.Range("N103").Value = -3010623.9
$ws.Range("H131").Value = 5292382
$ws.Range("I131").Value = 456.66666
$ws.Range("J131").Value = 5849427
$ws.Range("K131").Value = 1369.99998
$ws.Range("L131").Value = 17548281
$ws.Range("M131").Value = 3670.00002
$ws.Range("N131").Value = -17558361
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5280.4863
$ws.Range("I70").Value = 5243.375
$ws.Range("K70").Value = 5243.375
$ws.Range("M70").Value = -4973.375
$ws.Range("H73").Value = 5280.4863
$ws.Range("I73").Value = 5243.375
$ws.Range("K73").Value = 5243.375
$ws.Range("M73").Value = -4307.375
$ws.Range("H102").Value = 3151.7
$ws.Range("I102").Value = 3159.6428
$ws.Range("J102").Value = 3133.1667
$ws.Range("K102").Value = 3159.6428
$ws.Range("L102").Value = 3133.1667
$ws.Range("M102").Value = -1537.6428
$ws.Range("N102").Value = -6377.1667
$ws.Range("H113").Value = 1289.1333
$ws.Range("I113").Value = 1211.8334
$ws.Range("J113").Value = 1340.6666
$ws.Range("K113").Value = 1211.8334
$ws.Range("L113").Value = 1340.6666
$ws.Range("M113").Value = 958.1666
$ws.Range("N113").Value = -5680.6666
$ws.Range("H122").Value = 1788.1428
$ws.Range("I122").Value = 1623.4
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 4870.200000000001
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -2420.200000000001
$ws.Range("N122").Value = -11500
$ws.Range("H123").Value = 11192
$ws.Range("J123").Value = 11192
$ws.Range("L123").Value = 11192
$ws.Range("N123").Value = -16092
$ws.Range("H132").Value = 2309.6365
$ws.Range("I132").Value = 1709.3658
$ws.Range("J132").Value = 4067.5715
$ws.Range("K132").Value = 5128.097400000001
$ws.Range("L132").Value = 12202.7145
$ws.Range("M132").Value = -2598.097400000001
$ws.Range("N132").Value = -17262.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2879.6453
$ws.Range("I7").Value = 2577.889
$ws.Range("J7").Value = 3003.0908
$ws.Range("K7").Value = 2577.889
$ws.Range("L7").Value = 3003.0908
$ws.Range("M7").Value = -2465.889
$ws.Range("N7").Value = -3227.0908
$ws.Range("H40").Value = 3550
$ws.Range("I40").Value = 2900
$ws.Range("J40").Value = 3588.2354
$ws.Range("K40").Value = 2900
$ws.Range("L40").Value = 3588.2354
$ws.Range("M40").Value = -2764
$ws.Range("N40").Value = -3860.2354
$ws.Range("H61").Value = 4216.5454
$ws.Range("I61").Value = 4172.125
$ws.Range("J61").Value = 4335
$ws.Range("K61").Value = 4172.125
$ws.Range("L61").Value = 4335
$ws.Range("M61").Value = -3970.125
$ws.Range("N61").Value = -4739
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H113").Value = 4216.5454
$ws.Range("I113").Value = 4172.125
$ws.Range("J113").Value = 4335
$ws.Range("K113").Value = 4172.125
$ws.Range("L113").Value = 4335
$ws.Range("M113").Value = -2002.125
$ws.Range("N113").Value = -8675
$ws.Range("H122").Value = 3594.0386
$ws.Range("I122").Value = 2775
$ws.Range("K122").Value = 8325
$ws.Range("M122").Value = -5875
$ws.Range("H126").Value = 2879.6453
$ws.Range("I126").Value = 2577.889
$ws.Range("J126").Value = 3003.0908
$ws.Range("K126").Value = 7733.667
$ws.Range("L126").Value = 9009.2724
$ws.Range("M126").Value = -5263.667
$ws.Range("N126").Value = -13949.2724
$ws.Range("N108").ClearContents()
